# S-matrix_doc.xlsx update
#
# Commit message: "Added check for column dataset, added WorkingArea to sample
# parameters, fixed beamcurrent measurement bug, started a manual."
# Of that, this file's change is the "WorkingArea" sample-parameter addition
# (plus the resulting scroll/selection state left on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S-Matrix")

# New row describing the "WorkingArea" sample parameter, appended right after
# the existing "ColMode" parameter (row 44) in the "Properties of WF" table:
# # | Name | Meaning | Vartype | notes
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "WorkingArea"
$ws.Range("C45").Value = "Edges of working area (um)"
$ws.Range("D45").Value = "string"
$ws.Range("E45").Value = "bottomleft U, bottomleft V, upperright U, upperright V"

# Leave the sheet scrolled/selected where the author ended up after typing
# the new row's notes.
$ws.Range("E46").Select()
